# "error solve ifrs list" -- replace the placeholder/garbage financial figures
# for AJ네트웍스 with the corrected per-period figures, and drop the
# (erroneous, not-yet-available) 2020E/2021E rows down to just their label
# columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 : 2014/12 (IFRS연결) ------------------------------------------
$ws.Range("D2").Value  = 10112
$ws.Range("E2").Value  = 770
$ws.Range("F2").Value  = 770
$ws.Range("G2").Value  = 364
$ws.Range("H2").Value  = 281
$ws.Range("I2").Value  = 166
$ws.Range("J2").Value  = 115
$ws.Range("K2").Value  = 14379
$ws.Range("L2").Value  = 11546
$ws.Range("M2").Value  = 2833
$ws.Range("N2").Value  = 1616
$ws.Range("O2").Value  = 1216
$ws.Range("P2").Value  = 341
$ws.Range("Q2").Value  = -802
$ws.Range("R2").Value  = -321
$ws.Range("S2").Value  = 1352
$ws.Range("T2").Value  = 371
$ws.Range("U2").Value  = -1173
$ws.Range("V2").Value  = 10160
$ws.Range("W2").Value  = 7.61
$ws.Range("X2").Value  = 2.77
$ws.Range("Y2").Value  = 13.05
$ws.Range("Z2").Value  = 2.87
$ws.Range("AA2").Value = 407.6
$ws.Range("AB2").Value = 439.72
$ws.Range("AC2").Value = 546
# AD2 (PER) no longer applies for this period -- the cell is removed outright.
$ws.Range("AD2").ClearContents()
$ws.Range("AE2").Value = 4743
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = 0
# AH2 (현금배당수익률) likewise removed outright, not just zeroed.
$ws.Range("AH2").ClearContents()
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 34082240

# ---- Row 3 : 2015/12 (IFRS연결) ------------------------------------------
$ws.Range("D3").Value  = 10556
$ws.Range("E3").Value  = 743
$ws.Range("F3").Value  = 743
$ws.Range("G3").Value  = 383
$ws.Range("H3").Value  = 290
$ws.Range("I3").Value  = 188
$ws.Range("J3").Value  = 102
$ws.Range("K3").Value  = 17252
$ws.Range("L3").Value  = 13272
$ws.Range("M3").Value  = 3980
$ws.Range("N3").Value  = 2628
$ws.Range("O3").Value  = 1352
$ws.Range("P3").Value  = 468
$ws.Range("Q3").Value  = -1063
$ws.Range("R3").Value  = -744
$ws.Range("S3").Value  = 2312
$ws.Range("T3").Value  = 669
$ws.Range("U3").Value  = -1733
$ws.Range("V3").Value  = 11612
$ws.Range("W3").Value  = 7.03
$ws.Range("X3").Value  = 2.75
$ws.Range("Y3").Value  = 8.869999999999999
$ws.Range("Z3").Value  = 1.83
$ws.Range("AA3").Value = 333.46
$ws.Range("AB3").Value = 509.19
$ws.Range("AC3").Value = 477
$ws.Range("AD3").Value = 19.75
$ws.Range("AE3").Value = 5613
$ws.Range("AF3").Value = 1.68
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 46822295

# ---- Row 4 : 2016/12 (IFRS연결) ------------------------------------------
$ws.Range("D4").Value  = 12539
$ws.Range("E4").Value  = 619
$ws.Range("F4").Value  = 619
$ws.Range("G4").Value  = 286
$ws.Range("H4").Value  = 171
$ws.Range("I4").Value  = 135
$ws.Range("J4").Value  = 36
$ws.Range("K4").Value  = 20813
$ws.Range("L4").Value  = 16538
$ws.Range("M4").Value  = 4275
$ws.Range("N4").Value  = 2800
$ws.Range("O4").Value  = 1475
$ws.Range("P4").Value  = 468
$ws.Range("Q4").Value  = -1867
$ws.Range("R4").Value  = -851
$ws.Range("S4").Value  = 3202
$ws.Range("T4").Value  = 431
$ws.Range("U4").Value  = -2298
$ws.Range("V4").Value  = 14544
$ws.Range("W4").Value  = 4.93
$ws.Range("X4").Value  = 1.36
$ws.Range("Y4").Value  = 4.98
$ws.Range("Z4").Value  = 0.9
$ws.Range("AA4").Value = 386.87
$ws.Range("AB4").Value = 538.03
$ws.Range("AC4").Value = 288
$ws.Range("AD4").Value = 22.54
$ws.Range("AE4").Value = 5980
$ws.Range("AF4").Value = 1.09
$ws.Range("AG4").Value = 60
$ws.Range("AH4").Value = 0.92
$ws.Range("AI4").Value = 20.8
$ws.Range("AJ4").Value = 46822295

# ---- Row 5 : 2017/12 (IFRS연결) ------------------------------------------
$ws.Range("D5").Value  = 8439
$ws.Range("E5").Value  = 224
$ws.Range("F5").Value  = 224
$ws.Range("G5").Value  = 8
$ws.Range("H5").Value  = 206
$ws.Range("I5").Value  = 150
$ws.Range("J5").Value  = 56
$ws.Range("K5").Value  = 23542
$ws.Range("L5").Value  = 19071
$ws.Range("M5").Value  = 4472
$ws.Range("N5").Value  = 2888
$ws.Range("O5").Value  = 1584
$ws.Range("P5").Value  = 468
$ws.Range("Q5").Value  = -1485
$ws.Range("R5").Value  = -894
$ws.Range("S5").Value  = 1872
$ws.Range("T5").Value  = 527
$ws.Range("U5").Value  = -2012
$ws.Range("V5").Value  = 16367
$ws.Range("W5").Value  = 2.66
$ws.Range("X5").Value  = 2.45
$ws.Range("Y5").Value  = 5.27
$ws.Range("Z5").Value  = 0.93
$ws.Range("AA5").Value = 426.47
$ws.Range("AB5").Value = 564.05
$ws.Range("AC5").Value = 320
$ws.Range("AD5").Value = 21.76
$ws.Range("AE5").Value = 6168
$ws.Range("AF5").Value = 1.13
$ws.Range("AG5").Value = 86
$ws.Range("AH5").Value = 1.23
$ws.Range("AI5").Value = 26.85
$ws.Range("AJ5").Value = 46822295

# ---- Row 6 : 2018/12 (IFRS연결) ------------------------------------------
$ws.Range("D6").Value  = 10567
$ws.Range("E6").Value  = -213
$ws.Range("F6").Value  = -213
$ws.Range("G6").Value  = -495
$ws.Range("H6").Value  = 373
$ws.Range("I6").Value  = 11
$ws.Range("K6").Value  = 25751
$ws.Range("L6").Value  = 20833
$ws.Range("M6").Value  = 4918
$ws.Range("N6").Value  = 2871
$ws.Range("P6").Value  = 468
$ws.Range("Q6").Value  = -748
$ws.Range("R6").Value  = -552
$ws.Range("S6").Value  = 1275
$ws.Range("T6").Value  = 386
$ws.Range("U6").Value  = -1134
$ws.Range("V6").Value  = 9928
$ws.Range("W6").Value  = -2.01
$ws.Range("X6").Value  = 3.53
$ws.Range("Y6").Value  = 0.39
$ws.Range("Z6").Value  = 1.51
$ws.Range("AA6").Value = 423.61
$ws.Range("AB6").Value = 561.8200000000001
$ws.Range("AC6").Value = 24
$ws.Range("AD6").Value = 186.12
$ws.Range("AE6").Value = 6131
$ws.Range("AF6").Value = 0.73
$ws.Range("AG6").Value = 100
$ws.Range("AH6").Value = 2.22
$ws.Range("AI6").Value = 413.15
$ws.Range("AJ6").Value = 46822295

# ---- Rows 7-9 : 2019/12(E), 2020/12(E), 2021/12(E) -----------------------
# These estimate rows no longer have any figures -- only the A (index),
# B ("연간") and C (period label) cells survive; all of D:AJ are cleared.
$ws.Range("D7:AJ9").ClearContents()
